# Generate Report for handoff
#
# The handoff package for 27144277-3f91-4e5e-8751-84420dbded78.md is now
# ready, so the per-locale handoff status sheets (and the Overview roll-up)
# need to reflect that:
#   - Status changes from "Handoff transform failed" to "Ready for handoff"
#   - The newly produced .xlf handoff file is recorded (with a hyperlink)
#   - The "Latest Handoff Datetime" is stamped
#   - The "Handoff Reason" changes from "Ignored" to "Include"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/855c8f4d0ba0c921e6ccb138b8201dd75ac3341b"

# Update the overall status shown on the Overview roll-up sheet too, since
# it references the very same status text.
$wsOverview.Range("B2:C2").Replace("Handoff transform failed", "Ready for handoff") | Out-Null

$locales = @(
    @{ Sheet = $wsZhCn; File = "27144277-3f91-4e5e-8751-84420dbded78.7e22bf946c64c3552130ef0bbe6c3d7d3d7403c8.zh-cn.xlf"; Stamp = "2016-01-11 08:24:15" },
    @{ Sheet = $wsDeDe; File = "27144277-3f91-4e5e-8751-84420dbded78.7e22bf946c64c3552130ef0bbe6c3d7d3d7403c8.de-de.xlf"; Stamp = "2016-01-11 08:24:32" }
)

foreach ($locale in $locales) {
    $ws = $locale.Sheet

    # Status: handoff is now ready
    $ws.Range("B2").Replace("Handoff transform failed", "Ready for handoff") | Out-Null

    # Re-create the existing hyperlinks together with the new one so that
    # relationship ids stay ordered left-to-right / top-to-bottom
    # (A2, then the new C2 link, then A3).
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/27144277-3f91-4e5e-8751-84420dbded78.md", "", "", "27144277-3f91-4e5e-8751-84420dbded78.md") | Out-Null

    # Record the produced handoff (.xlf) file, with a hyperlink to it, like
    # the other tracked files on this sheet.
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/e2e/$($locale.File)", "", "", $locale.File) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/.localization-config", "", "", ".localization-config") | Out-Null

    # Latest Handoff Datetime
    $ws.Range("D2").Value = $locale.Stamp

    # Handoff Reason: no longer ignored, now included in the handoff
    $ws.Range("H2").Replace("Ignored", "Include") | Out-Null
}

Write-Host "Updated handoff status for Overview, zh-cn and de-de"
